$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 4,23
$arr[0,0] = 0.105952380952381
$arr[0,1] = 0.036734693877551
$arr[0,2] = 0.116156462585034
$arr[0,3] = 0.897448979591837
$arr[0,4] = 0.0204081632653061
$arr[0,5] = 0.951700680272109
$arr[0,6] = 0.937244897959184
$arr[0,7] = 0.210544217687075
$arr[0,8] = 0.949319727891156
$arr[0,9] = 0.0022108843537415
$arr[0,10] = 0.00578231292517007
$arr[0,11] = 0.00119047619047619
$arr[0,12] = 0.000170068027210884
$arr[0,13] = 0.0335034013605442
$arr[0,14] = 0.0518707482993197
$arr[0,15] = 0.0217687074829932
$arr[0,16] = 0.0256802721088435
$arr[0,17] = 0.998979591836735
$arr[0,18] = 0.968027210884354
$arr[0,19] = 0.122789115646259
$arr[0,20] = 0.00306122448979592
$arr[0,21] = 0.00289115646258503
$arr[0,22] = 0.000340136054421769
$arr[1,0] = 0.0192176870748299
$arr[1,1] = 0.870068027210884
$arr[1,2] = 0.000170068027210884
$arr[1,3] = 0.00255102040816327
$arr[1,4] = 0.0506802721088435
$arr[1,5] = 0.000340136054421769
$arr[1,6] = 0.0569727891156463
$arr[1,7] = 0.00714285714285714
$arr[1,8] = 0
$arr[1,9] = 0.000680272108843537
$arr[1,10] = 0
$arr[1,11] = 0.00952380952380952
$arr[1,12] = 0.000170068027210884
$arr[1,13] = 0.0022108843537415
$arr[1,14] = 0
$arr[1,15] = 0.00697278911564626
$arr[1,16] = 0.836054421768707
$arr[1,17] = 0.000510204081632653
$arr[1,18] = 0.000170068027210884
$arr[1,19] = 0.0105442176870748
$arr[1,20] = 0.0326530612244898
$arr[1,21] = 0.0481292517006803
$arr[1,22] = 0.0653061224489796
$arr[2,0] = 0.869897959183674
$arr[2,1] = 0.0564625850340136
$arr[2,2] = 0.054421768707483
$arr[2,3] = 0.0719387755102041
$arr[2,4] = 0.00986394557823129
$arr[2,5] = 0.0472789115646259
$arr[2,6] = 0.00561224489795918
$arr[2,7] = 0.779591836734694
$arr[2,8] = 0.0488095238095238
$arr[2,9] = 0.0569727891156463
$arr[2,10] = 0.992176870748299
$arr[2,11] = 0.988775510204082
$arr[2,12] = 0.990136054421769
$arr[2,13] = 0.0564625850340136
$arr[2,14] = 0.947789115646258
$arr[2,15] = 0.971258503401361
$arr[2,16] = 0.0542517006802721
$arr[2,17] = 0.000510204081632653
$arr[2,18] = 0.00153061224489796
$arr[2,19] = 0.866156462585034
$arr[2,20] = 0.0178571428571429
$arr[2,21] = 0.904081632653061
$arr[2,22] = 0.930272108843537
$arr[3,0] = 0.00493197278911565
$arr[3,1] = 0.036734693877551
$arr[3,2] = 0.829251700680272
$arr[3,3] = 0.0280612244897959
$arr[3,4] = 0.919047619047619
$arr[3,5] = 0.000680272108843537
$arr[3,6] = 0.000170068027210884
$arr[3,7] = 0.00238095238095238
$arr[3,8] = 0.00187074829931973
$arr[3,9] = 0.940136054421769
$arr[3,10] = 0.00204081632653061
$arr[3,11] = 0.000510204081632653
$arr[3,12] = 0.00952380952380952
$arr[3,13] = 0.90765306122449
$arr[3,14] = 0.000340136054421769
$arr[3,15] = 0
$arr[3,16] = 0.0840136054421769
$arr[3,17] = 0
$arr[3,18] = 0.0302721088435374
$arr[3,19] = 0.000340136054421769
$arr[3,20] = 0.946428571428571
$arr[3,21] = 0.0448979591836735
$arr[3,22] = 0.00408163265306122

$ws.Range("B2:X5").Value = $arr
